$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the B1 header from "TotalCpmI" to "MeanCpmI"
$ws.Range("B1").Value = "MeanCpmI"

# Update the active selection to B2
$ws.Range("B2").Select()
